$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - values are textual (not numeric) in the source data,
# so we format the cells as Text before assigning to avoid Excel auto-converting
# the numeric-looking strings into real numbers, then restore the default style.
$priceUpdates = @{
    "D2" = "242.57"
    "D3" = "23.01"
    "D4" = "5.418"
    "D5" = "0.05890"
    "D6" = "3.440"
    "D7" = "6.522"
    "D8" = "0.8106"
    "D9" = "0.9349"
    "D10" = "0.1418"
    "D11" = "0.07419"
    "D12" = "0.03308"
    "D13" = "0.03051"
    "D15" = "3.856"
    "D16" = "0.001564"
    "D17" = "0.04678"
    "D18" = "0.0005955"
    "D19" = "0.005861"
    "D20" = "0.001263"
    "D21" = "0.004886"
    "D26" = "0.1331"
    "D40" = "0.03947"
    "D42" = "0.1071"
    "D43" = "0.002573"
    "D44" = "0.009101"
    "D45" = "0.00005184"
    "D47" = "0.6705"
    "D48" = "0.002392"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# Column E (Volume label) text corrections
$ws.Range("E20").Value = "19BitKanKANBestin24h"
$ws.Range("E44").Value = "43LocalTradersLCT"

